# Update gh-pages output data (想去人数 / "want to go" counts) on the
# "展览" and "全部类型" sheets to match the newly generated data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 656
$ws1.Range("F4").Value = 256
$ws1.Range("F6").Value = 10073
$ws1.Range("F10").Value = 5769
$ws1.Range("F11").Value = 11
$ws1.Range("F12").Value = 99
$ws1.Range("F13").Value = 183
$ws1.Range("F22").Value = 21
$ws1.Range("F23").Value = 1533

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 656
$ws4.Range("F5").Value = 256
$ws4.Range("F7").Value = 10073
$ws4.Range("F11").Value = 5769
$ws4.Range("F12").Value = 11
$ws4.Range("F13").Value = 99
$ws4.Range("F14").Value = 183
$ws4.Range("F23").Value = 21
$ws4.Range("F24").Value = 1533
